# Add a new "FSAE_Achilles" sheet by duplicating the existing Sedan_HambaLG
# template sheet (keeps all styles / layout / column widths identical),
# then tweak the two cells that differ for the new vehicle instance.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("Sedan_HambaLG")

# Copy the sheet so it lands immediately after the source sheet; Excel
# auto-names it "Sedan_HambaLG (2)" and makes it the active sheet.
$src.Copy($null, $src)

$new = $wb.Worksheets.Item($src.Index + 1)
$new.Name = "FSAE_Achilles"

# Instance name (row "Instance") now matches the new sheet/template name.
$new.Range("H3").Value = "FSAE_Achilles"

# rWheelCutout changes for the new vehicle.
$new.Range("H6").Value = 0.25

$new.Activate()
